$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 53

$ws.Cells.Item($row, 1).Value = "2025-08-24 12:59:52 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-24 18:29:52 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""

$ws.Range("A52:H52").Copy()
$ws.Range("A53:H53").PasteSpecial(-4122)
$excel.CutCopyMode = $false
